$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")
$ws.Activate()

# Update price per unit (CHF) values in column F
$ws.Range("F3").Value = 11
$ws.Range("F5").Value = 71
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 2
$ws.Range("F21").Value = 2

# Update the active selection on the BoM sheet
$ws.Range("I17").Select()
